$wb = $excel.ActiveWorkbook

# --- Table S6 "TAD conservation" sheet: rename the conservation bucket
# labels from the old generic "Group 1".."Group 5" to the new descriptive
# names used throughout the rest of the workbook. This also causes the
# now-unused "Group N" shared-string entries to be dropped automatically
# and the five new strings to be appended to the shared-string table.
$ws6 = $wb.Worksheets.Item("Table S6 TAD conservation")

$ws6.Range("A3").Value = "Unique"
$ws6.Range("A4").Value = "Rare"
$ws6.Range("A5").Value = "Moderately Conserved"
$ws6.Range("A6").Value = "Highly Conserved"
$ws6.Range("A7").Value = "Core"

$ws6.Range("A9").Value = "Unique"
$ws6.Range("A10").Value = "Rare"
$ws6.Range("A11").Value = "Moderately Conserved"
$ws6.Range("A12").Value = "Highly Conserved"
$ws6.Range("A13").Value = "Core"

# --- Move the active/selected tab from Table S1 to Table S6, and update
# the remembered selection on Table S6 to the newly highlighted block of
# renamed labels.
$ws6.Activate()
$ws6.Range("A9:A13").Select()
